$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet had two extra numbered log-entry rows (rows 14-18, items 3-7)
# that were removed once the sheet was finalized; the totals / footer rows
# below shift up by five rows as a result.
# ---------------------------------------------------------------------------
$ws.Range("A14:H18").EntireRow.Delete()

# Row 10 (the "Backhoe per hour" rate line) was manually resized down.
$ws.Rows(10).RowHeight = 33

# The print area shrank along with the sheet (used to go to row 30, now 25).
$ws.PageSetup.PrintArea = '$A$1:$H$25'

# Update the on-screen scroll position / selection left after the edits.
[void]$ws.Range("A14:XFD18").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
